# "adding intern projects domains"
# Append 8 new domain rows (237-244) to Tabelle1, mirroring the existing
# A/B column pattern: column B holds the bare domain, column A holds a
# formula that prefixes it with "https://www." (like all prior rows) -
# except for the very last row, which uses a plain "https://" prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Domains that follow the usual "https://www." + domain pattern.
$domains = @(
    "amjis.com",
    "ilamecca.de",
    "kalender.amjis.com",
    "nivontec.de",
    "nordstern.gmbh",
    "nma-hamburg.de",
    "smarttech-elektro.de"
)

$startRow = 237
$row = $startRow
foreach ($domain in $domains) {
    $ws.Cells.Item($row, 2).Value = $domain
    $row = $row + 1
}
$endRow = $row - 1

# Fill A237:A243 in one shot so the engine groups them as a shared formula,
# matching how the rest of the column is built.
$ws.Range("A" + $startRow + ":A" + $endRow).Formula = '="https://www." & B' + $startRow

# Highlight the newly-added block (light/background fill), same as the
# "schrott" block above it used a purple fill.
$ws.Range("B" + $startRow + ":B" + $endRow).Interior.ThemeColor = 2

# Final row: a new project domain that isn't a "www." site, so it gets its
# own one-off formula and no fill.
$lastRow = $endRow + 1
$ws.Cells.Item($lastRow, 2).Value = "dr-elrafei.de"
$ws.Range("A" + $lastRow).Formula = '="https://" & B' + $lastRow

# Leave the cursor on the last entered cell, like the source workbook.
$null = $ws.Range("B" + $lastRow).Select()
